# Refresh the cryptos price/volume snapshot (GitHub Actions scrape update).
# Values that look numeric (e.g. "1.000", "0.4747") are written with a
# leading apostrophe so Excel keeps them as literal text, matching the
# original sheet's inline-string storage instead of coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.602.98'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '1.923.94'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''247.16'
$ws.Range('E5').Value = '  +2.84%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').Value = '''0.4747'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = '''0.2907'
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').Value = '''0.06841'
$ws.Range('E9').Value = '  +3.97%  '
$ws.Range('D10').Value = '''105.48'
$ws.Range('E10').Value = '  -1.57%  '
$ws.Range('D11').Value = '''18.43'
$ws.Range('E11').Value = '  -3.51%  '
$ws.Range('D12').Value = '1.927.57'
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').Value = '''0.07699'
$ws.Range('E13').Value = '  +1.29%  '
$ws.Range('D14').Value = '''5.357'
$ws.Range('E14').Value = '  +4.61%  '
$ws.Range('D15').Value = '''0.6719'
$ws.Range('E15').Value = '  +2.43%  '
$ws.Range('D16').Value = '''290.88'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '30.615.78'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = '''0.000007625'
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').Value = '''12.96'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').Value = '''5.554'
$ws.Range('E21').Value = '  +4.99%  '
$ws.Range('D22').Value = '2.178.10'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('D23').Value = '''1.006'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').Value = '''6.473'
$ws.Range('E24').Value = '  +2.98%  '
$ws.Range('D25').Value = '''9.549'
$ws.Range('E25').Value = '  +3.79%  '
$ws.Range('D26').Value = '''167.34'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').Value = '''21.17'
$ws.Range('E27').Value = '  +5.68%  '
$ws.Range('D28').Value = '''2.123'
$ws.Range('E28').Value = '  +5.07%  '
$ws.Range('D29').Value = '''0.1074'
$ws.Range('E29').Value = '  -3.48%  '
$ws.Range('D30').Value = '''1.400'
$ws.Range('E30').Value = '  +3.47%  '
$ws.Range('D31').Value = '''4.183'
$ws.Range('E31').Value = '  +2.60%  '
$ws.Range('D32').Value = '''4.053'
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('D33').Value = '''0.05032'
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').Value = '''0.7318'
$ws.Range('E34').Value = '  -1.06%  '
$ws.Range('D35').Value = '''1.146'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').Value = '''0.02074'
$ws.Range('E36').Value = '  +7.00%  '
$ws.Range('D37').Value = '''0.9997'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').Value = '''2.733'
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').Value = '''2.684'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').Value = '''112.18'
$ws.Range('E40').Value = '  +5.00%  '
$ws.Range('D41').Value = '''2.045'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('D42').Value = '''0.8729'
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('D43').Value = '''0.4420'
$ws.Range('E43').Value = '  +7.13%  '
$ws.Range('D44').Value = '''5.940'
$ws.Range('E44').Value = '  +2.19%  '
$ws.Range('D45').Value = '''1.001'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').Value = '''67.99'
$ws.Range('E46').Value = '  -2.53%  '
$ws.Range('D47').Value = '''7.308'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''9.380'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D49').Value = '''48.76'
$ws.Range('E49').Value = '  +16.62%  '
$ws.Range('E50').Value = '  +3.96%  '
$ws.Range('D51').Value = '''35.07'
